$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the date formatting from the last existing row (A36) down to the new row (A37)
$ws.Range("A36").Copy()
$ws.Range("A37").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new row of data: 12/18/2025 -> serial 46009, Error Count 4
$ws.Range("A37").Value = 46009
$ws.Range("B37").Value = 4

# Move the selection to the new last row, matching the saved workbook view
$ws.Range("A37:B37").Select()
